$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header swap: average_doctor / average_doctor_old ---
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"

# --- Row 4 ---
$ws.Range("E4").Value = 0.361
$ws.Range("F4").Value = 0.077
$ws.Range("G4").Value = 0.277
$ws.Range("N4").Value = 0.402
$ws.Range("O4").Value = 0.064
$ws.Range("P4").Value = 0.253
$ws.Range("W4").Value = 0.285
$ws.Range("AI4").Value = 0.203
$ws.Range("AJ4").Value = 0.06
$ws.Range("AK4").Value = 0.245
$ws.Range("AU4").Value = 0.13
$ws.Range("AV4").Value = 0.018
$ws.Range("AW4").Value = 0.134
$ws.Range("BA4").Value = 1.759
$ws.Range("BB4").Value = 0.154
$ws.Range("BC4").Value = 0.392
$ws.Range("BG4").Value = 0.671
$ws.Range("BH4").Value = 0.148
$ws.Range("BI4").Value = 0.385
$ws.Range("BM4").Value = 0.634
$ws.Range("BN4").Value = 0.091
$ws.Range("BO4").Value = 0.302
$ws.Range("BP4").Value = 0.586
$ws.Range("BQ4").Value = 0.608

# --- Row 5 ---
$ws.Range("E5").Value = 0.478
$ws.Range("F5").Value = 0.089
$ws.Range("G5").Value = 0.298
$ws.Range("N5").Value = 0.75
$ws.Range("O5").Value = 0.082
$ws.Range("P5").Value = 0.286
$ws.Range("W5").Value = 0.296
$ws.Range("X5").Value = 0.128
$ws.Range("Y5").Value = 0.358
$ws.Range("AI5").Value = 0.25
$ws.Range("AJ5").Value = 0.091
$ws.Range("AK5").Value = 0.302
$ws.Range("AU5").Value = 0.298
$ws.Range("AV5").Value = 0.097
$ws.Range("AW5").Value = 0.311
$ws.Range("BA5").Value = 1.367
$ws.Range("BB5").Value = 0.102
$ws.Range("BC5").Value = 0.32
$ws.Range("BG5").Value = 0.406
$ws.Range("BH5").Value = 0.051
$ws.Range("BI5").Value = 0.226
$ws.Range("BM5").Value = 0.586
$ws.Range("BN5").Value = 0.092
$ws.Range("BO5").Value = 0.303
$ws.Range("BP5").Value = 0.456
$ws.Range("BQ5").Value = 0.45

# --- Row 6 ---
$ws.Range("E6").Value = 0.411
$ws.Range("N6").Value = 0.523
$ws.Range("W6").Value = 0.29
$ws.Range("AI6").Value = 0.224
$ws.Range("AU6").Value = 0.181
$ws.Range("BA6").Value = 1.526
$ws.Range("BG6").Value = 0.506
$ws.Range("BM6").Value = 0.609
$ws.Range("BP6").Value = 0.509
$ws.Range("BQ6").Value = 0.514

# --- Row 7 ---
$ws.Range("E7").Value = 0.449
$ws.Range("N7").Value = 0.639
$ws.Range("W7").Value = 0.294
$ws.Range("AI7").Value = 0.239
$ws.Range("AU7").Value = 0.237
$ws.Range("BA7").Value = 1.425
$ws.Range("BG7").Value = 0.441
$ws.Range("BM7").Value = 0.595
$ws.Range("BP7").Value = 0.475
$ws.Range("BQ7").Value = 0.473

# --- Row 8 ---
$ws.Range("E8").Value = 0.496
$ws.Range("F8").Value = 0.12
$ws.Range("G8").Value = 0.347
$ws.Range("N8").Value = 0.757
$ws.Range("O8").Value = 0.063
$ws.Range("P8").Value = 0.251
$ws.Range("W8").Value = 0.286
$ws.Range("X8").Value = 0.119
$ws.Range("Y8").Value = 0.345
$ws.Range("AI8").Value = 0.233
$ws.Range("AJ8").Value = 0.097
$ws.Range("AK8").Value = 0.312
$ws.Range("AU8").Value = 0.217
$ws.Range("AV8").Value = 0.061
$ws.Range("AW8").Value = 0.246
$ws.Range("BA8").Value = 1.607
$ws.Range("BB8").Value = 0.142
$ws.Range("BC8").Value = 0.376
$ws.Range("BG8").Value = 0.523
$ws.Range("BH8").Value = 0.109
$ws.Range("BI8").Value = 0.33
$ws.Range("BM8").Value = 0.676
$ws.Range("BN8").Value = 0.082
$ws.Range("BO8").Value = 0.286
$ws.Range("BP8").Value = 0.536
$ws.Range("BQ8").Value = 0.5570000000000001

# --- Row 9 ---
$ws.Range("E9").Value = 0.421
$ws.Range("F9").Value = 0.244
$ws.Range("G9").Value = 0.494
$ws.Range("N9").Value = 0.658
$ws.Range("O9").Value = 0.225
$ws.Range("P9").Value = 0.474
$ws.Range("W9").Value = 0.158
$ws.Range("X9").Value = 0.133
$ws.Range("Y9").Value = 0.365
$ws.Range("AI9").Value = 0.132
$ws.Range("AJ9").Value = 0.114
$ws.Range("AK9").Value = 0.338
$ws.Range("BA9").Value = 1.499
$ws.Range("BB9").Value = 0.233
$ws.Range("BC9").Value = 0.482
$ws.Range("BG9").Value = 0.526
$ws.Range("BH9").Value = 0.249
$ws.Range("BI9").Value = 0.499
$ws.Range("BM9").Value = 0.605
$ws.Range("BN9").Value = 0.239
$ws.Range("BO9").Value = 0.489
$ws.Range("BP9").Value = 0.5
$ws.Range("BQ9").Value = 0.513

# --- Row 10 ---
$ws.Range("E10").Value = 0.553
$ws.Range("F10").Value = 0.247
$ws.Range("G10").Value = 0.497
$ws.Range("N10").Value = 0.868
$ws.Range("O10").Value = 0.114
$ws.Range("P10").Value = 0.338
$ws.Range("W10").Value = 0.342
$ws.Range("X10").Value = 0.225
$ws.Range("Y10").Value = 0.474
$ws.Range("AI10").Value = 0.263
$ws.Range("AJ10").Value = 0.194
$ws.Range("AK10").Value = 0.44
$ws.Range("AU10").Value = 0.211
$ws.Range("AV10").Value = 0.166
$ws.Range("AW10").Value = 0.408
$ws.Range("BA10").Value = 1.868
$ws.Range("BB10").Value = 0.247
$ws.Range("BC10").Value = 0.497
$ws.Range("BG10").Value = 0.579
$ws.Range("BH10").Value = 0.244
$ws.Range("BI10").Value = 0.494
$ws.Range("BM10").Value = 0.842
$ws.Range("BN10").Value = 0.133
$ws.Range("BO10").Value = 0.365
$ws.Range("BP10").Value = 0.623
$ws.Range("BQ10").Value = 0.667

# --- Row 11 ---
$ws.Range("E11").Value = 0.579
$ws.Range("F11").Value = 0.244
$ws.Range("G11").Value = 0.494
$ws.Range("N11").Value = 0.895
$ws.Range("O11").Value = 0.094
$ws.Range("P11").Value = 0.307
$ws.Range("W11").Value = 0.342
$ws.Range("X11").Value = 0.225
$ws.Range("Y11").Value = 0.474
$ws.Range("AI11").Value = 0.263
$ws.Range("AJ11").Value = 0.194
$ws.Range("AK11").Value = 0.44
$ws.Range("AU11").Value = 0.342
$ws.Range("AV11").Value = 0.225
$ws.Range("AW11").Value = 0.474
$ws.Range("BA11").Value = 1.868
$ws.Range("BB11").Value = 0.247
$ws.Range("BC11").Value = 0.497
$ws.Range("BG11").Value = 0.579
$ws.Range("BH11").Value = 0.244
$ws.Range("BI11").Value = 0.494
$ws.Range("BM11").Value = 0.842
$ws.Range("BN11").Value = 0.133
$ws.Range("BO11").Value = 0.365
$ws.Range("BP11").Value = 0.623
$ws.Range("BQ11").Value = 0.667

# --- Row 12 ---
$ws.Range("E12").Value = 1.455
$ws.Range("F12").Value = 0.702
$ws.Range("G12").Value = 0.838
$ws.Range("N12").Value = 1.514
$ws.Range("O12").Value = 1.107
$ws.Range("P12").Value = 1.052
$ws.Range("W12").Value = 1.846
$ws.Range("X12").Value = 0.746
$ws.Range("Y12").Value = 0.863
$ws.Range("AI12").Value = 1.9
$ws.Range("AJ12").Value = 0.89
$ws.Range("AU12").Value = 2.923
$ws.Range("AV12").Value = 1.456
$ws.Range("AW12").Value = 1.206
$ws.Range("BA12").Value = 3.790999999999999
$ws.Range("BB12").Value = 0.443
$ws.Range("BC12").Value = 0.666
$ws.Range("BG12").Value = 1.091
$ws.Range("BH12").Value = 0.083
$ws.Range("BI12").Value = 0.287
$ws.Range("BM12").Value = 1.406
$ws.Range("BN12").Value = 0.491
$ws.Range("BO12").Value = 0.701
$ws.Range("BP12").Value = 1.264
$ws.Range("BQ12").Value = 1.328

# --- Row 13 ---
$ws.Range("E13").Value = 1.75
$ws.Range("F13").Value = 0.899
$ws.Range("G13").Value = 0.948
$ws.Range("N13").Value = 2.276
$ws.Range("O13").Value = 0.957
$ws.Range("P13").Value = 0.978
$ws.Range("W13").Value = 1.09
$ws.Range("X13").Value = 0.186
$ws.Range("Y13").Value = 0.431
$ws.Range("AI13").Value = 1.39
$ws.Range("AJ13").Value = 0.419
$ws.Range("AK13").Value = 0.647
$ws.Range("AU13").Value = 2.45
$ws.Range("AV13").Value = 0.743
$ws.Range("AW13").Value = 0.862
$ws.Range("BA13").Value = 2.693
$ws.Range("BB13").Value = 0.295
$ws.Range("BC13").Value = 0.543
$ws.Range("BG13").Value = 0.671
$ws.Range("BH13").Value = 0.091
$ws.Range("BI13").Value = 0.302
$ws.Range("BM13").Value = 1.07
$ws.Range("BN13").Value = 0.378
$ws.Range("BO13").Value = 0.615
$ws.Range("BP13").Value = 0.898
$ws.Range("BQ13").Value = 0.8169999999999999
